$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set width for new column BS (71) to match the other date columns (12 chars)
$ws.Columns.Item(71).ColumnWidth = 11.17

$ws.Range("A1").Copy()
$ws.Range("BS1").PasteSpecial(-4122)
$ws.Range("BS1").Formula = '="2024/11/18"'
$ws.Range("BS1").Copy()
$ws.Range("BS1").PasteSpecial(-4163)

$ws.Range("D2").Copy()
$ws.Range("BS2").PasteSpecial(-4122)
$ws.Range("BS2").Value = 120.1

$ws.Range("A3").Copy()
$ws.Range("BS3").PasteSpecial(-4122)
$ws.Range("BS3").Value = 175.8

$ws.Range("A4").Copy()
$ws.Range("BS4").PasteSpecial(-4122)
$ws.Range("BS4").Value = 141.1

$ws.Range("A5").Copy()
$ws.Range("BS5").PasteSpecial(-4122)
$ws.Range("BS5").Value = 153.8

$ws.Range("A6").Copy()
$ws.Range("BS6").PasteSpecial(-4122)
$ws.Range("BS6").Value = 195.2

$ws.Range("B7").Copy()
$ws.Range("BS7").PasteSpecial(-4122)
$ws.Range("BS7").Value = 134.1

$ws.Range("B8").Copy()
$ws.Range("BS8").PasteSpecial(-4122)
$ws.Range("BS8").Value = 139.5

$ws.Range("A9").Copy()
$ws.Range("BS9").PasteSpecial(-4122)
$ws.Range("BS9").Value = 171.4

$ws.Range("A10").Copy()
$ws.Range("BS10").PasteSpecial(-4122)
$ws.Range("BS10").Value = 149.3

$ws.Range("A11").Copy()
$ws.Range("BS11").PasteSpecial(-4122)
$ws.Range("BS11").Value = 142.7

$ws.Range("A12").Copy()
$ws.Range("BS12").PasteSpecial(-4122)
$ws.Range("BS12").Value = 144.6

$ws.Range("A13").Copy()
$ws.Range("BS13").PasteSpecial(-4122)
$ws.Range("BS13").Value = 141

$ws.Range("A14").Copy()
$ws.Range("BS14").PasteSpecial(-4122)
$ws.Range("BS14").Value = 329.9

$ws.Range("H15").Copy()
$ws.Range("BS15").PasteSpecial(-4122)
$ws.Range("BS15").Value = 125.2

$ws.Range("A16").Copy()
$ws.Range("BS16").PasteSpecial(-4122)
$ws.Range("BS16").Value = 212.9

$ws.Range("A17").Copy()
$ws.Range("BS17").PasteSpecial(-4122)
$ws.Range("BS17").Value = 146.8

$ws.Range("A18").Copy()
$ws.Range("BS18").PasteSpecial(-4122)
$ws.Range("BS18").Value = 230.6

$ws.Range("A19").Copy()
$ws.Range("BS19").PasteSpecial(-4122)
$ws.Range("BS19").Value = 149.7

$ws.Range("H20").Copy()
$ws.Range("BS20").PasteSpecial(-4122)
$ws.Range("BS20").Value = 136.8

$ws.Range("A21").Copy()
$ws.Range("BS21").PasteSpecial(-4122)
$ws.Range("BS21").Value = 172

$ws.Range("A22").Copy()
$ws.Range("BS22").PasteSpecial(-4122)
$ws.Range("BS22").Value = 202

$ws.Range("A23").Copy()
$ws.Range("BS23").PasteSpecial(-4122)
$ws.Range("BS23").Value = 164.3

$ws.Range("K24").Copy()
$ws.Range("BS24").PasteSpecial(-4122)
$ws.Range("BS24").Value = 123.3

$ws.Range("A25").Copy()
$ws.Range("BS25").PasteSpecial(-4122)
$ws.Range("BS25").Value = 161.8

$ws.Range("A26").Copy()
$ws.Range("BS26").PasteSpecial(-4122)
$ws.Range("BS26").Value = 143.5

$ws.Range("A27").Copy()
$ws.Range("BS27").PasteSpecial(-4122)
$ws.Range("BS27").Value = 162.6

$ws.Range("J28").Copy()
$ws.Range("BS28").PasteSpecial(-4122)
$ws.Range("BS28").Value = 116.9

$ws.Range("A29").Copy()
$ws.Range("BS29").PasteSpecial(-4122)
$ws.Range("BS29").Value = 145.6

$ws.Range("A30").Copy()
$ws.Range("BS30").PasteSpecial(-4122)
$ws.Range("BS30").Value = 141

$ws.Range("A31").Copy()
$ws.Range("BS31").PasteSpecial(-4122)
$ws.Range("BS31").Value = 150.2

$ws.Range("A32").Copy()
$ws.Range("BS32").PasteSpecial(-4122)
$ws.Range("BS32").Value = 205.3

$ws.Range("A33").Copy()
$ws.Range("BS33").PasteSpecial(-4122)
$ws.Range("BS33").Value = 161.9

$ws.Range("C34").Copy()
$ws.Range("BS34").PasteSpecial(-4122)
$ws.Range("BS34").Value = 137.1

$ws.Range("A35").Copy()
$ws.Range("BS35").PasteSpecial(-4122)
$ws.Range("BS35").Value = 231.5

$ws.Range("I36").Copy()
$ws.Range("BS36").PasteSpecial(-4122)
$ws.Range("BS36").Value = 129

$ws.Range("A37").Copy()
$ws.Range("BS37").PasteSpecial(-4122)
$ws.Range("BS37").Value = 150

$ws.Range("A38").Copy()
$ws.Range("BS38").PasteSpecial(-4122)
$ws.Range("BS38").Value = 158.5

$ws.Range("A39").Copy()
$ws.Range("BS39").PasteSpecial(-4122)
$ws.Range("BS39").Value = 153.4

$ws.Range("A40").Copy()
$ws.Range("BS40").PasteSpecial(-4122)
$ws.Range("BS40").Value = 209.1

$ws.Range("A41").Copy()
$ws.Range("BS41").PasteSpecial(-4122)
$ws.Range("BS41").Value = 165.5

$ws.Range("A42").Copy()
$ws.Range("BS42").PasteSpecial(-4122)
$ws.Range("BS42").Value = 174.6

$ws.Range("A43").Copy()
$ws.Range("BS43").PasteSpecial(-4122)
$ws.Range("BS43").Value = 230.3

$ws.Range("A44").Copy()
$ws.Range("BS44").PasteSpecial(-4122)
$ws.Range("BS44").Value = 179.7

$ws.Range("AJ45").Copy()
$ws.Range("BS45").PasteSpecial(-4122)
$ws.Range("BS45").Value = 120.1

$ws.Range("A46").Copy()
$ws.Range("BS46").PasteSpecial(-4122)
$ws.Range("BS46").Value = 148.7

$ws.Range("A47").Copy()
$ws.Range("BS47").PasteSpecial(-4122)
$ws.Range("BS47").Value = 181

$ws.Range("A48").Copy()
$ws.Range("BS48").PasteSpecial(-4122)
$ws.Range("BS48").Value = 192.2

$ws.Range("A49").Copy()
$ws.Range("BS49").PasteSpecial(-4122)
$ws.Range("BS49").Value = 156.6

$ws.Range("A50").Copy()
$ws.Range("BS50").PasteSpecial(-4122)
$ws.Range("BS50").Value = 143.9

$ws.Range("A51").Copy()
$ws.Range("BS51").PasteSpecial(-4122)
$ws.Range("BS51").Value = 287.4

$ws.Range("A52").Copy()
$ws.Range("BS52").PasteSpecial(-4122)
$ws.Range("BS52").Value = 162.9

$ws.Range("A53").Copy()
$ws.Range("BS53").PasteSpecial(-4122)
$ws.Range("BS53").Value = 175.8

$excel.CutCopyMode = 0
